# Update the student ID in A7 and move the active selection to A7
# (mirrors manual proctoring edits: corrected ID + cursor left where the
# last edit was made).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A7").Value = 22201765
$ws.Range("A7").Select()
